# Update Sage scrape results
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row: insert a new "The" column at G1, shifting the existing
#    G1:Y1 header labels one column to the right (into H1:Z1).
# ---------------------------------------------------------------------------

# Capture the existing header labels in G1:Y1 (19 cells) before overwriting.
$oldHeaders = @()
for ($col = 7; $col -le 25; $col++) {
    $oldHeaders += $ws.Cells.Item(1, $col).Value2
}

# Give the brand-new last header cell (Z1) the same look (bold/border) as its
# neighbour before we populate it, since it did not exist before.
$ws.Range("Y1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Write the shifted header labels into H1:Z1.
for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $ws.Cells.Item(1, 8 + $i).Value = $oldHeaders[$i]
}

# Finally set the new G1 label.
$ws.Range("G1").Value = "The"

# ---------------------------------------------------------------------------
# 2. Data rows 2-11: replace Title/Authors/Year/DOI/Access Type with the new
#    scrape results, and append a new Z column (value 0) for every row.
# ---------------------------------------------------------------------------

$rows = @{
    2 = @("Adaptive Torque and Position Control for a Legged Robot Based on a Series Elastic Actuator", "Qiuguo Zhu, Yichao Mao, Rong Xiong, Jun Wu", "2016", "10.5772/62204", "Open Access")
    3 = @("Cyber-Flirting: Playing at Love on the Internet", "Monica Therese Whitty", "2003", "10.1177/0959354303013003003", "Restricted")
    4 = @("A novel ensemble learning approach for fault detection of sensor data in cyber-physical system", "Ramesh Sneka Nandhini, Ramanathan Lakshmanan", "2023", "10.3233/JIFS-235809", "Restricted")
    5 = @("Moving beyond the sanctuary paradigm: Canada must face up to the reality of a contested and dangerous space environment", "Patrick Perron", "2023", "10.1177/00207020231178394", "Restricted")
    6 = @("Prioritizing investment in military cyber capability using risk analysis", "Cayt Rowe, Hossein Seif Zadeh, Ivan L. Garanovich, Li Jiang, Daniel Bilusich, Rick Nunes-Vaz, Anthony Ween", "2019", "10.1177/1548512917707077", "Restricted")
    7 = @("Fighting in Cyberspace: Internet Access and the Substitutability of Cyber and Military Operations", "Nadiya Kostyuk, Erik Gartzke", "2024", "10.1177/00220027231160993", "Restricted")
    8 = @("Internet of Things, cybersecurity and governing wicked problems: learning from climate change governance", "Madeline Carr, Feja Lesniewska", "2020", "10.1177/0047117820948247", "Open Access")
    9 = @("Digital Assays Part II: Digital Protein and Cell Assays", "Amar S. Basu", "2017", "10.1177/2472630317705681", "Restricted")
    10 = @("Cyber scares and prophylactic policies: Crossnational evidence on the effect of cyberattacks on public support for surveillance", "Amelia C Arsenault, Sarah E Kreps, Keren LG Snider, Daphna Canetti", "2024", "10.1177/00223433241233960", "Restricted")
    11 = @("Towards a Chronology of Robotic Art", "Eduardo Kac", "2001", "10.1177/135485650100700109", "Restricted")
}

# Force the Year column to text first so the purely-numeric-looking values
# (e.g. "2016") are not silently converted into numbers by Excel.
$ws.Range("D2:D11").NumberFormat = "@"

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
}

# Restore the Year column's original (default) formatting/style now that the
# text values are locked in, so no stray style index is introduced.
$ws.Range("A2").Copy()
$ws.Range("D2:D11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 26).Value = 0
}
